$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save old row 7 values (C# native CLR / Windows)
$a7 = $ws.Cells.Item(7, 1).Value2
$b7 = $ws.Cells.Item(7, 2).Value2

# Shift rows 8-11 up to rows 7-10
for ($r = 8; $r -le 11; $r++) {
    $ws.Cells.Item($r - 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r - 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r - 1, 5).Value2 = $ws.Cells.Item($r, 5).Value2
}

# Place old row7 (C# native CLR on windows) content into row 11, with NEW value
$ws.Cells.Item(11, 1).Value2 = $a7
$ws.Cells.Item(11, 2).Value2 = $b7
$ws.Cells.Item(11, 5).Value2 = 35647
